$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.497.87'
$ws.Range('E2').Value = '  -0.73%  '
$ws.Range('D3').Value = '2.924.13'
$ws.Range('E3').Value = '  -2.34%  '
$ws.Range('E4').Value = '  -0.10%  '
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '375.86'
$c.ClearFormats()
$ws.Range('E5').Value = '  +6.45%  '
$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '104.08'
$c.ClearFormats()
$ws.Range('E6').Value = '  -2.15%  '
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.542'
$c.ClearFormats()
$ws.Range('E7').Value = '  -2.44%  '
$ws.Range('E8').Value = '  -0.21%  '
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.585'
$c.ClearFormats()
$ws.Range('E9').Value = '  -3.52%  '
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '36.81'
$c.ClearFormats()
$ws.Range('E10').Value = '  -2.48%  '
$ws.Range('E11').Value = '  -0.61%  '
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.0836'
$c.ClearFormats()
$ws.Range('E12').Value = '  -1.74%  '
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '18.31'
$c.ClearFormats()
$ws.Range('E13').Value = '  -3.04%  '
$ws.Range('D14').Value = '3.387.60'
$ws.Range('E14').Value = '  -2.62%  '
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '7.38'
$c.ClearFormats()
$ws.Range('E15').Value = '  -2.25%  '
$ws.Range('D16').Value = '2.922.94'
$ws.Range('E16').Value = '  -2.49%  '
$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '0.938'
$c.ClearFormats()
$ws.Range('E17').Value = '  -6.56%  '
$ws.Range('D18').Value = '51.410.44'
$ws.Range('E18').Value = '  -0.94%  '
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '3.40'
$c.ClearFormats()
$ws.Range('E19').Value = '  +1.08%  '
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '7.31'
$c.ClearFormats()
$ws.Range('E20').Value = '  -1.44%  '
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '12.96'
$c.ClearFormats()
$ws.Range('E21').Value = '  -3.78%  '
$ws.Range('D22').Value = '0.0₃0945'
$ws.Range('E22').Value = '  -2.07%  '
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '68.32'
$c.ClearFormats()
$ws.Range('E23').Value = '  -0.73%  '
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '261.54'
$c.ClearFormats()
$ws.Range('E24').Value = '  -0.18%  '
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.78'
$c.ClearFormats()
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '0.168'
$c.ClearFormats()
$ws.Range('E26').Value = '  -5.02%  '
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '4.12'
$c.ClearFormats()
$ws.Range('E27').Value = '  -5.23%  '
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('B29').Value = 'EthereumClassic'
$ws.Range('C29').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '25.79'
$c.ClearFormats()
$ws.Range('E29').Value = '  -3.77%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '7.32'
$c.ClearFormats()
$ws.Range('E30').Value = '  -1.00%  '
$ws.Range('B31').Value = 'RenderToken'
$ws.Range('C31').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '7.02'
$c.ClearFormats()
$ws.Range('E31').Value = '  +10.14%  '
$ws.Range('E32').Value = '  -5.45%  '
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '9.82'
$c.ClearFormats()
$ws.Range('E33').Value = '  -3.02%  '
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '51.76'
$c.ClearFormats()
$ws.Range('E34').Value = '  +1.22%  '
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '2.09'
$c.ClearFormats()
$ws.Range('E35').Value = '  -3.33%  '
$c = $ws.Range('D36')
$c.NumberFormat = '@'
$c.Value = '34.01'
$c.ClearFormats()
$ws.Range('E36').Value = '  -4.80%  '
$ws.Range('E37').Value = '  +0.33%  '
$ws.Range('E38').Value = '  -1.27%  '
$ws.Range('E39').Value = '  -8.05%  '
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '16.92'
$c.ClearFormats()
$ws.Range('E40').Value = '  -2.35%  '
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '2.59'
$c.ClearFormats()
$ws.Range('E41').Value = '  -7.33%  '
$ws.Range('E42').Value = '  -5.38%  '
$ws.Range('E43').Value = '  -1.82%  '
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '124.35'
$c.ClearFormats()
$ws.Range('E44').Value = '  -0.21%  '
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '21.81'
$c.ClearFormats()
$ws.Range('E45').Value = '  -5.83%  '
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '2.05'
$c.ClearFormats()
$ws.Range('E46').Value = '  -6.24%  '
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '0.271'
$c.ClearFormats()
$ws.Range('E47').Value = '  +12.86%  '
$ws.Range('D48').Value = '2.019.00'
$ws.Range('E48').Value = '  -4.80%  '
$ws.Range('E49').Value = '  -1.03%  '
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '3.17'
$c.ClearFormats()
$ws.Range('E50').Value = '  -3.66%  '
$ws.Range('D51').Value = '3.213.64'
$ws.Range('E51').Value = '  -2.77%  '
